# Update cryptocurrency price/volume data per source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.705.37'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = '2.522.33'
$ws.Range('E3').Value = '  +0.23%  '
$c = $ws.Range('D4')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $s
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '592.29'
$c.Style = $s
$ws.Range('E5').Value = '  +2.50%  '
$c = $ws.Range('D6')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '176.71'
$c.Style = $s
$ws.Range('E6').Value = '  +6.04%  '
$ws.Range('E7').Value = '  -0.06%  '
$c = $ws.Range('D8')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.530'
$c.Style = $s
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('D9').Value = '2.521.00'
$ws.Range('E9').Value = '  +0.17%  '
$c = $ws.Range('D10')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.141'
$c.Style = $s
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('E11').Value = '  +2.59%  '
$c = $ws.Range('D12')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.16'
$c.Style = $s
$ws.Range('E12').Value = '  +1.21%  '
$c = $ws.Range('D13')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.344'
$c.Style = $s
$ws.Range('E13').Value = '  -0.98%  '
$c = $ws.Range('D14')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '26.82'
$c.Style = $s
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '2.976.62'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range('D17').Value = '67.474.11'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '2.509.82'
$ws.Range('E18').Value = '  -0.10%  '
$c = $ws.Range('D19')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.02'
$c.Style = $s
$ws.Range('E19').Value = '  +4.62%  '
$c = $ws.Range('D20')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.44'
$c.Style = $s
$ws.Range('E20').Value = '  +1.49%  '
$c = $ws.Range('D21')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '360.33'
$c.Style = $s
$ws.Range('E21').Value = '  +4.39%  '
$c = $ws.Range('D22')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.20'
$c.Style = $s
$ws.Range('E22').Value = '  +0.28%  '
$c = $ws.Range('D23')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.64'
$c.Style = $s
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('E25').Value = '  -0.02%  '
$c = $ws.Range('D26')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '71.10'
$c.Style = $s
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.656.92'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D29')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = $s
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').Value = '0.0₃0990'
$ws.Range('E30').Value = '  +1.66%  '
$c = $ws.Range('D31')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '550.17'
$c.Style = $s
$ws.Range('E31').Value = '  +6.09%  '
$ws.Range('E32').Value = '  +2.58%  '
$c = $ws.Range('D33')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.35'
$c.Style = $s
$ws.Range('E33').Value = '  +3.41%  '
$ws.Range('E34').Value = '  +3.30%  '
$ws.Range('E35').Value = '  +0.13%  '
$c = $ws.Range('D36')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $s
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +2.12%  '
$c = $ws.Range('D38')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '154.49'
$c.Style = $s
$ws.Range('E38').Value = '  -0.79%  '
$c = $ws.Range('D39')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '18.75'
$c.Style = $s
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('E40').Value = '  +1.90%  '
$c = $ws.Range('D41')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.356'
$c.Style = $s
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('E42').Value = '  +3.26%  '
$c = $ws.Range('D43')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = $s
$ws.Range('E43').Value = '  +3.08%  '
$c = $ws.Range('D44')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.53'
$c.Style = $s
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0282'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D47')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.562'
$c.Style = $s
$ws.Range('E47').Value = '  +1.69%  '
$c = $ws.Range('D48')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '147.03'
$c.Style = $s
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('E51').Value = '  +0.80%  '
